$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "clean the execute data": the sheet tracked both the expected result (col E)
# and the actual executed result / pass-fail columns (F "实际结果", G "是否通过").
# The per-row executed/actual data is stale test-run output, so clear it out for
# every data row while leaving the column headers (row 1) in place.
$ws.Range("F2:G18").ClearContents()

# Move the active selection like the author's last save did.
$null = $ws.Range("I14").Select()
